$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.862.45"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.945.97"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Formula = "'552.84"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Formula = "'133.19"
$ws.Range("E6").Value = "  +9.14%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Formula = "'0.512"
$ws.Range("E8").Value = "  +4.50%  "
$ws.Range("D9").Value = "2.939.55"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Formula = "'0.448"
$ws.Range("E12").Value = "  +4.14%  "
$ws.Range("E13").Value = "  +4.53%  "
$ws.Range("E14").Value = "  +4.92%  "
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "3.431.76"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("E17").Value = "  +8.13%  "
$ws.Range("D18").Value = "2.944.79"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").Value = "57.880.49"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Formula = "'417.27"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Formula = "'13.35"
$ws.Range("E21").Value = "  +4.31%  "
$ws.Range("D22").Formula = "'0.698"
$ws.Range("E22").Value = "  +7.39%  "
$ws.Range("D23").Formula = "'13.45"
$ws.Range("E23").Value = "  +7.35%  "
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("D25").Formula = "'78.96"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Formula = "'2.01"
$ws.Range("E29").Value = "  +5.13%  "
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("D31").Formula = "'25.49"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").Formula = "'5.95"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Formula = "'0.0970"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").Formula = "'5.69"
$ws.Range("E34").Value = "  +6.08%  "
$ws.Range("D35").Formula = "'0.946"
$ws.Range("E35").Value = "  +5.94%  "
$ws.Range("D36").Formula = "'2.07"
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("D37").Value = "0.0₃0702"
$ws.Range("E37").Value = "  +13.17%  "
$ws.Range("D38").Formula = "'48.37"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").Formula = "'8.78"
$ws.Range("E39").Value = "  +5.03%  "
$ws.Range("D40").Formula = "'2.69"
$ws.Range("E40").Value = "  +14.50%  "
$ws.Range("D41").Formula = "'381.34"
$ws.Range("E41").Value = "  +6.67%  "
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("D43").Formula = "'0.0347"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").Value = "2.701.28"
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Formula = "'124.07"
$ws.Range("E46").Value = "  +5.54%  "
$ws.Range("E47").Value = "  +3.85%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Formula = "'0.108"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("E50").Value = "  +2.29%  "
$ws.Range("D51").Formula = "'2.00"
$ws.Range("E51").Value = "  +3.32%  "
